$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.877.65"
$ws.Range("E2").Value = '  -1.18%  '
$ws.Range("D3").Value = "'2.336.85"
$ws.Range("E3").Value = '  +0.08%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = "'303.51"
$ws.Range("E5").Value = '  +0.23%  '
$ws.Range("D6").Value = "'93.85"
$ws.Range("E6").Value = '  -4.42%  '
$ws.Range("E7").Value = '  -1.41%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E9").Value = '  -1.90%  '
$ws.Range("D10").Value = "'34.01"
$ws.Range("E10").Value = '  -4.90%  '
$ws.Range("E11").Value = '  -2.18%  '
$ws.Range("D12").Value = "'18.69"
$ws.Range("E12").Value = '  -3.98%  '
$ws.Range("D13").Value = "'0.120"
$ws.Range("E13").Value = '  +1.03%  '
$ws.Range("D14").Value = "'6.71"
$ws.Range("E14").Value = '  -3.03%  '
$ws.Range("D15").Value = "'2.701.33"
$ws.Range("E15").Value = '  +0.27%  '
$ws.Range("D16").Value = "'2.320.36"
$ws.Range("E16").Value = '  -0.70%  '
$ws.Range("D17").Value = "'0.792"
$ws.Range("E17").Value = '  -0.16%  '
$ws.Range("D18").Value = "'42.817.59"
$ws.Range("E18").Value = '  -1.21%  '
$ws.Range("D19").Value = "'12.04"
$ws.Range("E19").Value = '  -5.99%  '
$ws.Range("E20").Value = '  +1.91%  '
$ws.Range("E21").Value = '  -1.61%  '
$ws.Range("D22").Value = "'67.83"
$ws.Range("E22").Value = '  -0.34%  '
$ws.Range("D23").Value = "'235.47"
$ws.Range("E23").Value = '  -0.97%  '
$ws.Range("E24").Value = '  -1.64%  '
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("E26").Value = '  -1.62%  '
$ws.Range("D27").Value = "'24.57"
$ws.Range("E27").Value = '  -2.03%  '
$ws.Range("E28").Value = '  -6.56%  '
$ws.Range("E29").Value = '  -0.18%  '
$ws.Range("D30").Value = "'31.27"
$ws.Range("E30").Value = '  -6.34%  '
$ws.Range("E31").Value = '  -0.04%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = "'4.97"
$ws.Range("E32").Value = '  -1.05%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = "'0.0739"
$ws.Range("E33").Value = '  +5.03%  '
$ws.Range("D34").Value = "'17.21"
$ws.Range("E34").Value = '  -3.92%  '
$ws.Range("D35").Value = "'4.40"
$ws.Range("E35").Value = '  -2.02%  '
$ws.Range("E36").Value = '  -1.05%  '
$ws.Range("E37").Value = '  +1.79%  '
$ws.Range("D38").Value = "'124.42"
$ws.Range("E38").Value = '  -24.39%  '
$ws.Range("E39").Value = '  -0.44%  '
$ws.Range("E40").Value = '  -1.05%  '
$ws.Range("D41").Value = "'22.32"
$ws.Range("E41").Value = '  +21.71%  '
$ws.Range("E42").Value = '  -1.75%  '
$ws.Range("D43").Value = "'1.934.93"
$ws.Range("E43").Value = '  -2.84%  '
$ws.Range("D44").Value = "'0.0282"
$ws.Range("E44").Value = '  -0.40%  '
$ws.Range("E45").Value = '  -5.22%  '
$ws.Range("E46").Value = '  +0.87%  '
$ws.Range("D47").Value = "'2.71"
$ws.Range("E47").Value = '  -3.54%  '
$ws.Range("E48").Value = '  -0.70%  '
$ws.Range("D49").Value = "'2.567.22"
$ws.Range("E49").Value = '  +0.14%  '
$ws.Range("D50").Value = "'52.72"
$ws.Range("E50").Value = '  -2.58%  '
$ws.Range("D51").Value = "'71.46"
$ws.Range("E51").Value = '  -1.99%  '
